$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

foreach ($col in @(10, 11)) {
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, $col)
        if ($cell.Value2 -eq "JV") {
            $cell.Value = "Junior Varsity"
        }
    }
}
